# Generate Report for Archive
#
# 1) The "Status" value used to read "Ready for handoff" wherever a file was
#    still waiting to be handed off; update it to "In Translation" so the
#    archived report reflects translation work that is now underway.
# 2) Narrow the "Status"-column width back down (it had been widened for the
#    previous handoff review pass) on every sheet that shows it: the
#    "Overview" sheet's zh-cn/de-de status columns, plus column C ("Status")
#    on each language sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Update status text -------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            # NOTE: compare with the literal on the LEFT -- Excel's COM
            # Range.Value can come back as a non-string (e.g. a boolean),
            # and "-eq" coerces the right-hand side to the left operand's
            # type, which would otherwise produce false positives.
            if ("Ready for handoff" -eq $cell.Value2) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2) Narrow the Status columns ------------------------------------------
# Target stored column width is ~13.41 "characters". The host's ColumnWidth
# setter snaps to a 1/6-character pixel grid (MDW-based), so feed it the
# character width (12.5) whose snapped result lands on the closest
# representable width to the target.
$newStatusColWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth   # column E: zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth   # column F: de-de

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth       # column C: Status

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth       # column C: Status
